# feat : playerInput, NaturePoint
# Insert a new "Hand" equip item (ID 1000) as the first data row, ahead of
# the existing 1001 (Hoe) row, and normalise the DropPercent (M) column for
# the pre-existing Equip-type rows (1001-1005) from -1 to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a fresh row at row 2; everything below (old rows 2-46) shifts
#    down to rows 3-47, carrying its formatting/row-height with it - this
#    matches a real Excel "Insert Sheet Rows" operation.
$ws.Rows("2:2").Insert()

# 2) Populate the new row 2 with the new "Hand" base-equip item.
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = "손"
$ws.Range("D2").Value = "기본 상태"
$ws.Range("E2").Value = "Equip"
$ws.Range("F2").Value = "Hand"
$ws.Range("G2").Value = -1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "Equip"
$ws.Range("L2").Value = "Hand"
$ws.Range("M2").Value = 0

# Match the look of the other "Item"-style rows (e.g. the stone/ore rows)
# by reusing the same font (fontId 1 - "Malgun Gothic", family 3) that the
# header's DropPercent cell (M1) already uses, for the new row's text +
# DropPercent cells.
$ws.Range("C2:F2").Font.Name = "맑은 고딕"
$ws.Range("C2:F2").Font.Family = 3
$ws.Range("K2:M2").Font.Name = "맑은 고딕"
$ws.Range("K2:M2").Font.Family = 3

# 3) The five pre-existing Equip rows (Hoe/Water/PickAxe/Axe/Sword, now at
#    rows 3-7 after the insert) get their DropPercent normalised from -1 to
#    0, with the same styling as the new Hand row's M cell.
foreach ($r in 3..7) {
    $ws.Range("M$r").Value = 0
    $ws.Range("M$r").Font.Name = "맑은 고딕"
    $ws.Range("M$r").Font.Family = 3
}

# 4) Leave the cursor where the author's saved session had it.
$ws.Range("O8").Select()
